# Auto-generated edit script applying the Sagittarius_Profits market-data refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2340.6365
$ws.Range("J2").Value = 4125.6
$ws.Range("L2").Value = 4125.6
$ws.Range("N2").Value = -4351.6
$ws.Range("H111").Value = 5761.615
$ws.Range("I111").Value = 6191.8
$ws.Range("J111").Value = 4327.6665
$ws.Range("K111").Value = 18575.4
$ws.Range("L111").Value = 12982.9995
$ws.Range("M111").Value = -15508.4
$ws.Range("N111").Value = -19116.9995
$ws.Range("H113").Value = 2839.8
$ws.Range("I113").Value = 1850
$ws.Range("J113").Value = 3499.6667
$ws.Range("K113").Value = 1850
$ws.Range("L113").Value = 3499.6667
$ws.Range("M113").Value = 1404
$ws.Range("N113").Value = -10007.6667
$ws.Range("H125").Value = 2995
$ws.Range("J125").Value = 2995
$ws.Range("L125").Value = 26955
$ws.Range("N125").Value = -31875
$ws.Range("H132").Value = 2094.5557
$ws.Range("I132").Value = 2168.875
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 6506.625
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -3976.625
$ws.Range("N132").Value = -9560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2045.9
$ws.Range("I2").Value = 1003
$ws.Range("J2").Value = 3610.25
$ws.Range("K2").Value = 1003
$ws.Range("L2").Value = 3610.25
$ws.Range("M2").Value = -890
$ws.Range("N2").Value = -3836.25
$ws.Range("H4").Value = 522.6
$ws.Range("I4").Value = 522.6
$ws.Range("K4").Value = 522.6
$ws.Range("M4").Value = -406.6
$ws.Range("H45").Value = 2808.9644
$ws.Range("I45").Value = 1937.6
$ws.Range("K45").Value = 1937.6
$ws.Range("M45").Value = -1560.6
$ws.Range("H98").Value = 5999
$ws.Range("J98").Value = 5999
$ws.Range("L98").Value = 5999
$ws.Range("N98").Value = -11989
$ws.Range("H110").Value = 1284.2
$ws.Range("I110").Value = 1314.3334
$ws.Range("K110").Value = 1314.3334
$ws.Range("M110").Value = 730.6666
$ws.Range("H116").Value = 2045.9
$ws.Range("I116").Value = 1003
$ws.Range("J116").Value = 3610.25
$ws.Range("K116").Value = 1003
$ws.Range("L116").Value = 3610.25
$ws.Range("M116").Value = 1291
$ws.Range("N116").Value = -8198.25
$ws.Range("H122").Value = 11864.917
$ws.Range("J122").Value = 2159.3333
$ws.Range("L122").Value = 6477.999899999999
$ws.Range("N122").Value = -11377.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2045.9
$ws.Range("I3").Value = 1003
$ws.Range("J3").Value = 3610.25
$ws.Range("K3").Value = 1003
$ws.Range("L3").Value = 3610.25
$ws.Range("M3").Value = -889
$ws.Range("N3").Value = -3838.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1401.875
$ws.Range("I35").Value = 1480.1666
$ws.Range("J35").Value = 1167
$ws.Range("K35").Value = 1480.1666
$ws.Range("L35").Value = 1167
$ws.Range("M35").Value = -1186.1666
$ws.Range("N35").Value = -1755
$ws.Range("H132").Value = 3184.611
$ws.Range("I132").Value = 3194.75
$ws.Range("K132").Value = 9584.25
$ws.Range("M132").Value = -7054.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2622.6047
$ws.Range("J68").Value = 2880.946
$ws.Range("L68").Value = 8642.838
$ws.Range("N68").Value = -10264.838
$ws.Range("H71").Value = 2622.6047
$ws.Range("J71").Value = 2880.946
$ws.Range("L71").Value = 25928.514
$ws.Range("N71").Value = -34040.514
$ws.Range("H81").Value = 2266.3333
$ws.Range("H84").Value = 2266.3333
$ws.Range("H92").Value = 1625
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H98").Value = 1400.8572
$ws.Range("I98").Value = 1600
$ws.Range("K98").Value = 4800
$ws.Range("M98").Value = -3302

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1928.6666
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1578.7858
$ws.Range("I122").Value = 1400.3334
$ws.Range("K122").Value = 4201.0002
$ws.Range("M122").Value = -1751.0002
$ws.Range("H126").Value = 4334
$ws.Range("J126").Value = 4295.3335
$ws.Range("L126").Value = 12886.0005
$ws.Range("N126").Value = -17826.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 332.42856
$ws.Range("I16").Value = 300
$ws.Range("J16").Value = 413.5
$ws.Range("K16").Value = 300
$ws.Range("L16").Value = 413.5
$ws.Range("M16").Value = -130
$ws.Range("N16").Value = -753.5
$ws.Range("H22").Value = 1798.1428
$ws.Range("I22").Value = 1695.6666
$ws.Range("J22").Value = 1875
$ws.Range("K22").Value = 1695.6666
$ws.Range("L22").Value = 1875
$ws.Range("M22").Value = -1400.6666
$ws.Range("N22").Value = -2465
$ws.Range("H27").Value = 1798.1428
$ws.Range("I27").Value = 1695.6666
$ws.Range("J27").Value = 1875
$ws.Range("K27").Value = 1695.6666
$ws.Range("L27").Value = 1875
$ws.Range("M27").Value = -1588.6666
$ws.Range("N27").Value = -2089
$ws.Range("H40").Value = 2573.2727
$ws.Range("J40").Value = 3407.25
$ws.Range("L40").Value = 3407.25
$ws.Range("N40").Value = -3679.25
$ws.Range("H55").Value = 876.8570999999999
$ws.Range("I55").Value = 690.46155
$ws.Range("J55").Value = 1179.75
$ws.Range("K55").Value = 690.46155
$ws.Range("L55").Value = 1179.75
$ws.Range("M55").Value = -517.46155
$ws.Range("N55").Value = -1525.75
$ws.Range("H93").Value = 1499.6666
$ws.Range("I93").Value = 1899.5
$ws.Range("K93").Value = 1899.5
$ws.Range("M93").Value = -651.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10355.667
$ws.Range("I62").Value = 10343
$ws.Range("J62").Value = 10400
$ws.Range("K62").Value = 10343
$ws.Range("L62").Value = 10400
$ws.Range("M62").Value = -9719
$ws.Range("N62").Value = -11648
$ws.Range("H65").Value = 10355.667
$ws.Range("I65").Value = 10343
$ws.Range("J65").Value = 10400
$ws.Range("K65").Value = 51715
$ws.Range("L65").Value = 52000
$ws.Range("M65").Value = -48595
$ws.Range("N65").Value = -58240
$ws.Range("H93").Value = 47416.668
$ws.Range("J93").Value = 47250
$ws.Range("L93").Value = 47250
$ws.Range("N93").Value = -52242
$ws.Range("H122").Value = 3574.5
$ws.Range("I122").Value = 1922.3334
$ws.Range("K122").Value = 5767.0002
$ws.Range("M122").Value = -3317.0002
$ws.Range("H128").Value = 72500
$ws.Range("J128").Value = 72500
$ws.Range("L128").Value = 72500
$ws.Range("N128").Value = -82460
